$wb = $excel.ActiveWorkbook

# The existing hyperlinks on all three sheets point at the same target URL;
# only the displayed text needs to track the new file name.
$linkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/056959abc5757d7b92ca261e34541f6b902966f1/e2e/046900f1-850e-4532-b4ee-2954445431cd.md"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "743b15eb-139e-47e3-840b-579365822d91.md"
$wsOverview.Range("B2").Value = "e2e\743b15eb-139e-47e3-840b-579365822d91.md"
$wsOverview.Range("G2").Value = "2016-08-18 11:00:52"

# Refresh the hyperlink display text on B2 (keep same target URL, just the
# shown text needs to track the new file name). Replace in place so the
# underlying relationship/address is preserved rather than duplicated.
$wsOverview.Range("B2").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "e2e\743b15eb-139e-47e3-840b-579365822d91.md")

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "743b15eb-139e-47e3-840b-579365822d91.md"
$wsZhCn.Range("G2").Value = "743b15eb-139e-47e3-840b-579365822d91.b2d28a441e03d7be4e2c767a5fb99be0ce9734a1.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-18 11:00:48"

$wsZhCn.Range("A2").Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "743b15eb-139e-47e3-840b-579365822d91.md")

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "743b15eb-139e-47e3-840b-579365822d91.md"
$wsDeDe.Range("G2").Value = "743b15eb-139e-47e3-840b-579365822d91.b2d28a441e03d7be4e2c767a5fb99be0ce9734a1.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-18 11:00:52"

$wsDeDe.Range("A2").Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "743b15eb-139e-47e3-840b-579365822d91.md")
